$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-62 down to 51-63.
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the weekly "Arveja Verde" market entry.
$ws.Cells.Item(50, 1).Value2 = 10
$ws.Cells.Item(50, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value2 = "La Araucanía"
$ws.Cells.Item(50, 4).Value2 = 44522
$ws.Cells.Item(50, 5).Value2 = 9
$ws.Cells.Item(50, 6).Value2 = 100112022
$ws.Cells.Item(50, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(50, 8).Value2 = "Sin especificar"
$ws.Cells.Item(50, 9).Value2 = "Primera"
$ws.Cells.Item(50, 10).Value2 = 80
$ws.Cells.Item(50, 11).Value2 = 16000
$ws.Cells.Item(50, 12).Value2 = 16000
$ws.Cells.Item(50, 13).Value2 = 16000
$ws.Cells.Item(50, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(50, 16).Value2 = 640
$ws.Cells.Item(50, 17).Value2 = 25
$ws.Cells.Item(50, 18).Value2 = "Hortaliza"
